$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert a new data row at the top of the data block (row 16),
# which pushes the existing rows 16-54 down to 17-55 (preserving all of their
# data/formatting - matches the diff's row-shift pattern exactly), then fill
# the newly inserted row 16 with this week's new record.
$ws.Rows(16).Insert()

$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 45177
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 100112026
$ws.Cells.Item(16, 7).Value = "Haba"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 400
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 13).Value = 10500
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 420
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
